$d = $word.ActiveDocument

# --- Part 1: turn "III ... article" into "III ... article et vue :" by
#     appending a new run (same size formatting) after the existing run,
#     inside the very same paragraph. ---
$p13 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Contains("article")) {
        $p13 = $cand
    }
}

$r13 = $p13.Range
$r13.MoveEnd(1, -1)
$insertAt = $r13.End
$r13.InsertAfter(" et vue :")

# The engine merges same-formatted appended text into the prior run; force the
# appended text to materialize as its own <w:r> (as in the target) by
# toggling a character-level property on just the new span, then reverting
# it so the visible formatting stays identical to the original heading run.
$newRun = $d.Range($insertAt, $insertAt + 9)
$newRun.Bold = 1
$newRun.Bold = 0

# --- Part 2: drop the "IV - GITHUB :" run entirely, leaving its paragraph
#     (and paragraph formatting) empty but present. ---
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$rLast.MoveEnd(1, -1)
$rLast.Delete()
